# Apply the edits to the "OceanData" sheet's server/IP table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OceanData")
$ws.Activate() | Out-Null

# New gateway rows 3 and 4 ("R" / "Gateway R") - previously only had B/C filled in.
$ws.Range("D3").Value = "Gateway R"
$ws.Range("A3").Value = "R"
$ws.Range("D4").Value = "Gateway R"
$ws.Range("A4").Value = "R"

# Row 8 used to be "FS02" (Fileshare); it is replaced by the old row-9 SQL01/MSSQL entry.
$ws.Range("A8").Value = "SQL01"
$ws.Range("D8").Value = "MSSQL"

# Row 9 becomes SQL02.
$ws.Range("A9").Value = "SQL02"

# Row 10 becomes the old BU01/Linux/Backup entry.
$ws.Range("A10").Value = "BU01"
$ws.Range("B10").Value = "Linux"
$ws.Range("D10").Value = "Backup"

# Row 11 becomes the old WEB01/Apache entry.
$ws.Range("A11").Value = "WEB01"
$ws.Range("D11").Value = "Apache"

# Row 12 becomes the "-"/"-" MSSQL listener IP entry (192.168.20.20), replacing WEB01.
$ws.Range("A12").Value = "-"
$ws.Range("B12").Value = "-"
$ws.Range("C12").Value = "192.168.20.20"
$ws.Range("D12").Value = "MSSQL listener IP"
$ws.Range("A12:B12").HorizontalAlignment = -4108  # xlCenter, matches the "-"/"-" placeholder rows

# Remove the old trailing rows 13 and 14 (Fileshare listener IP / 192.168.20.30 MSSQL listener IP),
# shrinking the table down to A1:D12.
$ws.Rows.Item(14).Delete() | Out-Null
$ws.Rows.Item(13).Delete() | Out-Null

# Leave the selection where the author ended up (cell E12, just past the table).
$ws.Range("E12").Select() | Out-Null

Write-Host "edit complete"
